# Apply updated cryptocurrency price/volume data to the worksheet.
# Values in column D occasionally look like plain numbers (e.g. "7.50", "1.00")
# so we use Excel's leading-apostrophe quote-prefix to force them to remain
# text, matching the original inlineStr cell contents exactly (no rounding,
# no loss of trailing zeros, no reinterpretation of "." separated values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.861.46"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "'3.387.64"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'571.90"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "'141.57"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'3.385.68"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'7.50"
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "'0.392"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "'3.969.99"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "'28.18"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D17").Value = "'3.390.28"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "'60.970.41"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").Value = "'14.15"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("D21").Value = "'9.07"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").Value = "'387.87"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'0.560"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'73.13"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "'3.533.29"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'0.179"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'7.37"
$ws.Range("E30").Value = "  -5.61%  "
$ws.Range("D31").Value = "'8.13"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  -6.94%  "
$ws.Range("D33").Value = "'2.15"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'23.80"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'6.94"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").Value = "'3.418.33"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").Value = "'167.28"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "'5.06"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "'1.52"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "'0.0780"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "'26.89"
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("D43").Value = "'0.786"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D46").Value = "'41.85"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'1.69"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").Value = "'2.532.43"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("E49").Value = "  -3.91%  "
$ws.Range("D50").Value = "'6.86"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "'22.88"
$ws.Range("E51").Value = "  -2.44%  "
